$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# Two runs in the document end with a trailing space ("Test link before
# bookmark<nbsp>: " and "Test bookmark<nbsp>: ") but were written out
# without 'xml:space="preserve"', so the trailing space silently gets
# dropped when the part is reloaded/reserialized. Re-touch each affected
# run (toggle a character property around a self Find/Replace of its own
# text) so it gets re-materialized with xml:space="preserve" on its
# <w:t>; the visible text and formatting are left exactly as they were.
function Fix-TrailingSpaceRun($searchText, $fullLength) {
    $found = $d.Content
    $ok = $found.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
    if (-not $ok) {
        return
    }
    $start = $found.Start

    $runRange = $d.Range($start, $start + $fullLength)
    $runRange.Bold = 1
    $runRange.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                            $true, 1, $false, $searchText, 2)

    $restoreRange = $d.Range($start, $start + $fullLength)
    $restoreRange.Bold = 0
}

Fix-TrailingSpaceRun "Test link before bookmark" (("Test link before bookmark" + $nbsp + ": ").Length)
Fix-TrailingSpaceRun "Test bookmark" (("Test bookmark" + $nbsp + ": ").Length)
